$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bloco MEC-3B (quinta/coluna E e terça/coluna B -> terça/coluna C)
$ws.Range("E2").Value = "-"
$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("C4").Value = "MEC-3B-Ens. Dest. não Dest."
$ws.Range("E4").Value = "-"
$ws.Range("C6").Value = "MEC-3B-Ens. Dest. não Dest."
$ws.Range("C7").Value = "MEC-3B-Ens. Dest. não Dest."
$ws.Range("C8").Value = "MEC-3B-Ens. Dest. não Dest."

# Bloco MEC-3A (quarta/coluna D -> segunda/coluna B)
$ws.Range("B11").Value = "MEC-3A-Ens. Dest. não Dest."
$ws.Range("D11").Value = "-"
$ws.Range("B12").Value = "MEC-3A-Ens. Dest. não Dest."
$ws.Range("D12").Value = "-"
$ws.Range("B14").Value = "MEC-3A-Ens. Dest. não Dest."
$ws.Range("B15").Value = "MEC-3A-Ens. Dest. não Dest."
$ws.Range("D15").Value = "-"
$ws.Range("D16").Value = "-"
